$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ E=3; G=9.156959333333335; H=27.470878; I=0.969469463764299; J=0.9694694637642989; K=3; M=13.267299; N=39.801897; O=0.4248136128385448; P=0.4248136128385448; Q=121.488117406174; R=1093.393056655566; S=0.4118438254383586; T=0.4118438254383585 }
    3 = @{ E=3; G=9.156959333333335; H=27.470878; I=0.969469463764299; J=0.9694694637642989; K=3; M=10.340832; N=31.022496; O=0.3311093088108164; P=0.3311093088108164; Q=94.69057809683203; R=852.215202871488; S=0.3210003640601898; T=0.3210003640601898 }
    4 = @{ E=3; G=9.156959333333335; H=27.470878; I=0.969469463764299; J=0.9694694637642989; K=3; M=7.622739666666667; N=22.868219; O=0.2440770783506388; P=0.2440770783506388; Q=69.80111713625357; R=628.2100542262821; S=0.2366252742657506; T=0.2366252742657506 }
    5 = @{ E=3; G=0.288371; H=0.865113; I=0.03053053623570109; J=0.03053053623570109; K=3; M=13.267299; N=39.801897; O=0.4248136128385448; P=0.4248136128385448; Q=3.825904279929; R=34.43313851936099; S=0.01296978740018629; T=0.01296978740018629 }
    6 = @{ E=3; G=0.288371; H=0.865113; I=0.03053053623570109; J=0.03053053623570109; K=3; M=10.340832; N=31.022496; O=0.3311093088108164; P=0.3311093088108164; Q=2.981996064672; R=26.837964582048; S=0.01010894475062657; T=0.01010894475062657 }
    7 = @{ E=3; G=0.288371; H=0.865113; I=0.03053053623570109; J=0.03053053623570109; K=3; M=7.622739666666667; N=22.868219; O=0.2440770783506388; P=0.2440770783506388; Q=2.198177060416333; R=19.783593543747; S=0.007451804084888231; T=0.007451804084888231 }
}

foreach ($rowNum in $data.Keys) {
    $rowData = $data[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
